# Add season-record columns (Wins / Losses / Ties) to the sheet.
#
# The sheet currently ends at column AC (dimension A1:AC50). We append
# three new columns: AD = Wins, AE = Losses, AF = Ties, populated for
# every data row (2-50) with the team's season record (89-73-0).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastDataRow = 50

# --- Header row -------------------------------------------------------
# Copy the formatting of the existing last header cell (AC1, which
# carries the bold/bordered/centered header style) onto the three new
# header cells so they pick up the same style index instead of Excel
# fabricating a brand-new one, then set their text.
$ws.Range("AC1").Copy($ws.Range("AD1"))
$ws.Range("AC1").Copy($ws.Range("AE1"))
$ws.Range("AC1").Copy($ws.Range("AF1"))

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# --- Data rows ---------------------------------------------------------
# Season record is the same for every player row in this sheet: 89 wins,
# 73 losses, 0 ties.
for ($row = 2; $row -le $lastDataRow; $row++) {
    $ws.Cells.Item($row, 30).Value = 89   # AD
    $ws.Cells.Item($row, 31).Value = 73   # AE
    $ws.Cells.Item($row, 32).Value = 0    # AF
}
